# docs: Added event_status endpoint in excel
#
# Adds a new row (16) to Sheet1 describing the "GET /events_status" API,
# following the same visual pattern as the "miners" service block above it
# (row 15) but highlighted with a yellow fill instead of red.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content ------------------------------------------------------
# SERVICE column repeats "miners" (same service as row 15), ENDPOINT is the
# new route, and Res documents the JSON response shape. Req (D16) stays
# blank, matching the source row's layout.
$ws.Range("B16").Value = "miners"
$ws.Range("C16").Value = "GET /events_status"
$jsonText = "[`n  {`n    EndpointId: INT (PK),`n    CompanyId: INT (FK),`n    status: STR`n  }`n]"
$ws.Range("E16").Value = $jsonText

# --- Formatting -------------------------------------------------------------
# SERVICE/ENDPOINT cells: Arial font, left/top aligned, yellow fill (no wrap)
$rngText = $ws.Range("B16:C16")
$rngText.Font.Name = "Arial"
$rngText.Interior.Color = 65535   # RGB(255,255,0) / FFFF00 yellow
$rngText.HorizontalAlignment = -4131  # xlLeft
$rngText.VerticalAlignment = -4160    # xlTop

# Req/Res cells: default font, left/top aligned, yellow fill, wrap text
$rngWrap = $ws.Range("D16:E16")
$rngWrap.Interior.Color = 65535
$rngWrap.HorizontalAlignment = -4131
$rngWrap.VerticalAlignment = -4160
$rngWrap.WrapText = $true

# Taller row to fit the wrapped JSON body
$ws.Rows(16).RowHeight = 100.8

# Column C needs to widen slightly to fit "GET /events_status"
$ws.Columns("C").ColumnWidth = 18.053385416666668

# --- Viewport / selection ---------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2
[void]$ws.Range("D11").Select()
